$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3403.6
$ws.Range("J113").Value = 3004.5
$ws.Range("L113").Value = 3004.5
$ws.Range("N113").Value = -9512.5

$ws.Range("H137").Value = 2198.8965
$ws.Range("I137").Value = 2203.2593
$ws.Range("K137").Value = 6609.777900000001
$ws.Range("M137").Value = -4059.777900000001

$ws.Range("H138").Value = 2042.2593
$ws.Range("I138").Value = 2112.6667
$ws.Range("J138").Value = 2012.614
$ws.Range("K138").Value = 6338.000100000001
$ws.Range("L138").Value = 6037.842000000001
$ws.Range("M138").Value = -1198.000100000001
$ws.Range("N138").Value = -16317.842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 596567.5600000001
$ws.Range("I32").Value = 714790.8
$ws.Range("J32").Value = 14545.385
$ws.Range("K32").Value = 714790.8
$ws.Range("L32").Value = 14545.385
$ws.Range("M32").Value = -714503.8
$ws.Range("N32").Value = -15119.385

$ws.Range("H45").Value = 4922.1113
$ws.Range("I45").Value = 5037.375
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 5037.375
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -4660.375
$ws.Range("N45").Value = -4754

$ws.Range("H131").Value = 39700
$ws.Range("J131").Value = 39700
$ws.Range("L131").Value = 39700
$ws.Range("N131").Value = -49780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 55780
$ws.Range("J126").Value = 55780
$ws.Range("L126").Value = 55780
$ws.Range("N126").Value = -65660

$ws.Range("H132").Value = 42746.43
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 42746.43
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 42746.43
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -52866.43

$ws.Range("H134").Value = 3137.0278
$ws.Range("I134").Value = 3098.3684
$ws.Range("J134").Value = 3180.2354
$ws.Range("K134").Value = 9295.1052
$ws.Range("L134").Value = 9540.706200000001
$ws.Range("M134").Value = -6760.1052
$ws.Range("N134").Value = -14610.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4721.1187
$ws.Range("I31").Value = 1082.1177
$ws.Range("J31").Value = 9670.16
$ws.Range("K31").Value = 1082.1177
$ws.Range("L31").Value = 9670.16
$ws.Range("M31").Value = -787.1177
$ws.Range("N31").Value = -10260.16

$ws.Range("H34").Value = 4721.1187
$ws.Range("I34").Value = 1082.1177
$ws.Range("J34").Value = 9670.16
$ws.Range("K34").Value = 1082.1177
$ws.Range("L34").Value = 9670.16
$ws.Range("M34").Value = -880.1177
$ws.Range("N34").Value = -10074.16

$ws.Range("H52").Value = 59793.332
$ws.Range("J52").Value = 59793.332
$ws.Range("L52").Value = 59793.332
$ws.Range("N52").Value = -60381.332

$ws.Range("H60").Value = 9302.429
$ws.Range("J60").Value = 10670.667
$ws.Range("L60").Value = 10670.667
$ws.Range("N60").Value = -11692.667

$ws.Range("H111").Value = 41111
$ws.Range("J111").Value = 41111
$ws.Range("L111").Value = 41111
$ws.Range("N111").Value = -49291

$ws.Range("H123").Value = 38749.5
$ws.Range("J123").Value = 38749.5
$ws.Range("L123").Value = 38749.5
$ws.Range("N123").Value = -48549.5

$ws.Range("H132").Value = 2820.2
$ws.Range("I132").Value = 2625.75
$ws.Range("J132").Value = 3042.4285
$ws.Range("K132").Value = 7877.25
$ws.Range("L132").Value = 9127.2855
$ws.Range("M132").Value = -5347.25
$ws.Range("N132").Value = -14187.2855

$ws.Range("H133").Value = 35666.668
$ws.Range("J133").Value = 35666.668
$ws.Range("L133").Value = 35666.668
$ws.Range("N133").Value = -40726.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6088.3335
$ws.Range("I3").Value = 5306
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 15918
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -15806
$ws.Range("N3").Value = -30224

$ws.Range("H98").Value = 484
$ws.Range("I98").Value = 396.6
$ws.Range("J98").Value = 629.6667
$ws.Range("K98").Value = 1189.8
$ws.Range("L98").Value = 1889.0001
$ws.Range("M98").Value = 308.1999999999998
$ws.Range("N98").Value = -4885.0001

$ws.Range("H123").Value = 5604
$ws.Range("J123").Value = 9966
$ws.Range("L123").Value = 29898
$ws.Range("N123").Value = -34798

$ws.Range("H131").Value = 4208.75
$ws.Range("I131").Value = 442.72726
$ws.Range("J131").Value = 5865.8
$ws.Range("K131").Value = 1328.18178
$ws.Range("L131").Value = 17597.4
$ws.Range("M131").Value = 3711.81822
$ws.Range("N131").Value = -27677.4

$ws.Range("H133").Value = 12275.556
$ws.Range("I133").Value = 1500
$ws.Range("J133").Value = 13622.5
$ws.Range("K133").Value = 4500
$ws.Range("L133").Value = 40867.5
$ws.Range("M133").Value = 560
$ws.Range("N133").Value = -50987.5

$ws.Range("H136").Value = 3428.4285
$ws.Range("J136").Value = 3977.7778
$ws.Range("L136").Value = 11933.3334
$ws.Range("N136").Value = -22133.3334

$ws.Range("H139").Value = 2548.6978
$ws.Range("I139").Value = 1367.4762
$ws.Range("J139").Value = 3676.2273
$ws.Range("K139").Value = 4102.4286
$ws.Range("L139").Value = 11028.6819
$ws.Range("M139").Value = 1037.5714
$ws.Range("N139").Value = -21308.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 926.6667
$ws.Range("I122").Value = 940
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 2820
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -370
$ws.Range("N122").Value = -7600

$ws.Range("H132").Value = 3813.2222
$ws.Range("I132").Value = 3454.2964
$ws.Range("J132").Value = 4890
$ws.Range("K132").Value = 10362.8892
$ws.Range("L132").Value = 14670
$ws.Range("M132").Value = -7832.889200000001
$ws.Range("N132").Value = -19730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2609
$ws.Range("I122").Value = 2563
$ws.Range("J122").Value = 2827.5
$ws.Range("K122").Value = 7689
$ws.Range("L122").Value = 8482.5
$ws.Range("M122").Value = -5239
$ws.Range("N122").Value = -13382.5

$ws.Range("H124").Value = 46955.285
$ws.Range("J124").Value = 46955.285
$ws.Range("L124").Value = 46955.285
$ws.Range("N124").Value = -56775.285

$ws.Range("H127").Value = 63607.25
$ws.Range("J127").Value = 63607.25
$ws.Range("L127").Value = 63607.25
$ws.Range("N127").Value = -73527.25

$ws.Range("H132").Value = 1653.56
$ws.Range("I132").Value = 1621.6061
$ws.Range("J132").Value = 1715.5883
$ws.Range("K132").Value = 4864.8183
$ws.Range("L132").Value = 5146.7649
$ws.Range("M132").Value = -2334.8183
$ws.Range("N132").Value = -10206.7649

$ws.Range("H135").Value = 104949.69
$ws.Range("J135").Value = 104949.69
$ws.Range("L135").Value = 104949.69
$ws.Range("N135").Value = -115089.69

$ws.Range("H136").Value = 1882.0618
$ws.Range("I136").Value = 1598.9219
$ws.Range("J136").Value = 2948
$ws.Range("K136").Value = 4796.7657
$ws.Range("L136").Value = 8844
$ws.Range("M136").Value = -2246.7657
$ws.Range("N136").Value = -13944
